$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Prescaler" label in A4 and default prescaler value 1 in B4
$ws.Range("A4").Value = "Prescaler"
$ws.Range("B4").Value = 1

# Take the prescaler into account in the F_Timer calculation
$ws.Range("B5").Formula = "=B2/B4/2^B3"

# Update the current selection to match the author's edit location
$ws.Range("B4").Select()
